$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 12.150543
$ws.Range("H2").Value = 36.451629
$ws.Range("I2").Value = 0.05382648024663175
$ws.Range("J2").Value = 0.05382648024663175
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 24.43393833333333
$ws.Range("N2").Value = 73.301815
$ws.Range("O2").Value = 0.1197493722400791
$ws.Range("P2").Value = 0.1197493722400791
$ws.Range("Q2").Value = 296.885618378515
$ws.Range("R2").Value = 2671.970565406636
$ws.Range("S2").Value = 0.006445687219427171
$ws.Range("T2").Value = 0.006445687219427172

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 12.150543
$ws.Range("H3").Value = 36.451629
$ws.Range("I3").Value = 0.05382648024663175
$ws.Range("J3").Value = 0.05382648024663175
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 86.43758666666668
$ws.Range("N3").Value = 259.31276
$ws.Range("O3").Value = 0.4236258027695808
$ws.Range("P3").Value = 0.4236258027695808
$ws.Range("Q3").Value = 1050.26361360956
$ws.Range("R3").Value = 9452.372522486043
$ws.Range("S3").Value = 0.02280228590474036
$ws.Range("T3").Value = 0.02280228590474036

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 12.150543
$ws.Range("H4").Value = 36.451629
$ws.Range("I4").Value = 0.05382648024663175
$ws.Range("J4").Value = 0.05382648024663175
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 33.80250733333333
$ws.Range("N4").Value = 101.407522
$ws.Range("O4").Value = 0.1656642076314483
$ws.Range("P4").Value = 0.1656642076314483
$ws.Range("Q4").Value = 410.718818861482
$ws.Range("R4").Value = 3696.469369753338
$ws.Range("S4").Value = 0.008917121199648053
$ws.Range("T4").Value = 0.008917121199648053

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 12.150543
$ws.Range("H5").Value = 36.451629
$ws.Range("I5").Value = 0.05382648024663175
$ws.Range("J5").Value = 0.05382648024663175
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 25.01395866666667
$ws.Range("N5").Value = 75.041876
$ws.Range("O5").Value = 0.1225920196207674
$ws.Range("P5").Value = 0.1225920196207674
$ws.Range("Q5").Value = 303.933180379556
$ws.Range("R5").Value = 2735.398623416004
$ws.Range("S5").Value = 0.006598696922511926
$ws.Range("T5").Value = 0.006598696922511927

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 12.150543
$ws.Range("H6").Value = 36.451629
$ws.Range("I6").Value = 0.05382648024663175
$ws.Range("J6").Value = 0.05382648024663175
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 25.369122
$ws.Range("N6").Value = 76.107366
$ws.Range("O6").Value = 0.1243326553557499
$ws.Range("P6").Value = 0.1243326553557499
$ws.Range("Q6").Value = 308.248607733246
$ws.Range("R6").Value = 2774.237469599214
$ws.Range("S6").Value = 0.006692389217517548
$ws.Range("T6").Value = 0.006692389217517548

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 12.150543
$ws.Range("H7").Value = 36.451629
$ws.Range("I7").Value = 0.05382648024663175
$ws.Range("J7").Value = 0.05382648024663175
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 8.985195333333333
$ws.Range("N7").Value = 26.955586
$ws.Range("O7").Value = 0.04403594238237437
$ws.Range("P7").Value = 0.04403594238237437
$ws.Range("Q7").Value = 109.175002261066
$ws.Range("R7").Value = 982.5750203495941
$ws.Range("S7").Value = 0.002370299782786688
$ws.Range("T7").Value = 0.002370299782786688

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 22.297341
$ws.Range("H8").Value = 66.89202300000001
$ws.Range("I8").Value = 0.09877644027011076
$ws.Range("J8").Value = 0.09877644027011076
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 24.43393833333333
$ws.Range("N8").Value = 73.301815
$ws.Range("O8").Value = 0.1197493722400791
$ws.Range("P8").Value = 0.1197493722400791
$ws.Range("Q8").Value = 544.811854991305
$ws.Range("R8").Value = 4903.306694921746
$ws.Range("S8").Value = 0.01182841671445543
$ws.Range("T8").Value = 0.01182841671445544

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 22.297341
$ws.Range("H9").Value = 66.89202300000001
$ws.Range("I9").Value = 0.09877644027011076
$ws.Range("J9").Value = 0.09877644027011076
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 86.43758666666668
$ws.Range("N9").Value = 259.31276
$ws.Range("O9").Value = 0.4236258027695808
$ws.Range("P9").Value = 0.4236258027695808
$ws.Range("Q9").Value = 1927.328345123721
$ws.Range("R9").Value = 17345.95510611348
$ws.Range("S9").Value = 0.04184424880414722
$ws.Range("T9").Value = 0.04184424880414722

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 22.297341
$ws.Range("H10").Value = 66.89202300000001
$ws.Range("I10").Value = 0.09877644027011076
$ws.Range("J10").Value = 0.09877644027011076
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 33.80250733333333
$ws.Range("N10").Value = 101.407522
$ws.Range("O10").Value = 0.1656642076314483
$ws.Range("P10").Value = 0.1656642076314483
$ws.Range("Q10").Value = 753.7060326663341
$ws.Range("R10").Value = 6783.354293997007
$ws.Range("S10").Value = 0.01636372071000298
$ws.Range("T10").Value = 0.01636372071000298

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 22.297341
$ws.Range("H11").Value = 66.89202300000001
$ws.Range("I11").Value = 0.09877644027011076
$ws.Range("J11").Value = 0.09877644027011076
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 25.01395866666667
$ws.Range("N11").Value = 75.041876
$ws.Range("O11").Value = 0.1225920196207674
$ws.Range("P11").Value = 0.1225920196207674
$ws.Range("Q11").Value = 557.7447661505721
$ws.Range("R11").Value = 5019.702895355149
$ws.Range("S11").Value = 0.01210920330366297
$ws.Range("T11").Value = 0.01210920330366298

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 22.297341
$ws.Range("H12").Value = 66.89202300000001
$ws.Range("I12").Value = 0.09877644027011076
$ws.Range("J12").Value = 0.09877644027011076
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 25.369122
$ws.Range("N12").Value = 76.107366
$ws.Range("O12").Value = 0.1243326553557499
$ws.Range("P12").Value = 0.1243326553557499
$ws.Range("Q12").Value = 565.6639641046021
$ws.Range("R12").Value = 5090.975676941419
$ws.Range("S12").Value = 0.0122811371053715
$ws.Range("T12").Value = 0.0122811371053715

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 22.297341
$ws.Range("H13").Value = 66.89202300000001
$ws.Range("I13").Value = 0.09877644027011076
$ws.Range("J13").Value = 0.09877644027011076
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 8.985195333333333
$ws.Range("N13").Value = 26.955586
$ws.Range("O13").Value = 0.04403594238237437
$ws.Range("P13").Value = 0.04403594238237437
$ws.Range("Q13").Value = 200.345964298942
$ws.Range("R13").Value = 1803.113678690478
$ws.Range("S13").Value = 0.004349713632470641
$ws.Range("T13").Value = 0.004349713632470641

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 54.98715833333333
$ws.Range("H14").Value = 164.961475
$ws.Range("I14").Value = 0.2435911869821439
$ws.Range("J14").Value = 0.2435911869821439
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 24.43393833333333
$ws.Range("N14").Value = 73.301815
$ws.Range("O14").Value = 0.1197493722400791
$ws.Range("P14").Value = 0.1197493722400791
$ws.Range("Q14").Value = 1343.552835841903
$ws.Range("R14").Value = 12091.97552257713
$ws.Range("S14").Value = 0.02916989172432746
$ws.Range("T14").Value = 0.02916989172432746

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 54.98715833333333
$ws.Range("H15").Value = 164.961475
$ws.Range("I15").Value = 0.2435911869821439
$ws.Range("J15").Value = 0.2435911869821439
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 86.43758666666668
$ws.Range("N15").Value = 259.31276
$ws.Range("O15").Value = 0.4236258027695808
$ws.Range("P15").Value = 0.4236258027695808
$ws.Range("Q15").Value = 4752.957263991223
$ws.Range("R15").Value = 42776.61537592101
$ws.Range("S15").Value = 0.1031915121329058
$ws.Range("T15").Value = 0.1031915121329058

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 54.98715833333333
$ws.Range("H16").Value = 164.961475
$ws.Range("I16").Value = 0.2435911869821439
$ws.Range("J16").Value = 0.2435911869821439
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 33.80250733333333
$ws.Range("N16").Value = 101.407522
$ws.Range("O16").Value = 0.1656642076314483
$ws.Range("P16").Value = 0.1656642076314483
$ws.Range("Q16").Value = 1858.703822801661
$ws.Range("R16").Value = 16728.33440521495
$ws.Range("S16").Value = 0.04035434097740082
$ws.Range("T16").Value = 0.04035434097740083

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 54.98715833333333
$ws.Range("H17").Value = 164.961475
$ws.Range("I17").Value = 0.2435911869821439
$ws.Range("J17").Value = 0.2435911869821439
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 25.01395866666667
$ws.Range("N17").Value = 75.041876
$ws.Range("O17").Value = 0.1225920196207674
$ws.Range("P17").Value = 0.1225920196207674
$ws.Range("Q17").Value = 1375.446505747456
$ws.Range("R17").Value = 12379.0185517271
$ws.Range("S17").Value = 0.02986233557396099
$ws.Range("T17").Value = 0.02986233557396099

$ws.Range("E18").Value = 3
$ws.Range("G18").Value = 54.98715833333333
$ws.Range("H18").Value = 164.961475
$ws.Range("I18").Value = 0.2435911869821439
$ws.Range("J18").Value = 0.2435911869821439
$ws.Range("K18").Value = 3
$ws.Range("M18").Value = 25.369122
$ws.Range("N18").Value = 76.107366
$ws.Range("O18").Value = 0.1243326553557499
$ws.Range("P18").Value = 0.1243326553557499
$ws.Range("Q18").Value = 1394.97592819165
$ws.Range("R18").Value = 12554.78335372485
$ws.Range("S18").Value = 0.03028633909874894
$ws.Range("T18").Value = 0.03028633909874894

$ws.Range("E19").Value = 3
$ws.Range("G19").Value = 54.98715833333333
$ws.Range("H19").Value = 164.961475
$ws.Range("I19").Value = 0.2435911869821439
$ws.Range("J19").Value = 0.2435911869821439
$ws.Range("K19").Value = 3
$ws.Range("M19").Value = 8.985195333333333
$ws.Range("N19").Value = 26.955586
$ws.Range("O19").Value = 0.04403594238237437
$ws.Range("P19").Value = 0.04403594238237437
$ws.Range("Q19").Value = 494.0703584499278
$ws.Range("R19").Value = 4446.63322604935
$ws.Range("S19").Value = 0.01072676747479987
$ws.Range("T19").Value = 0.01072676747479987

$ws.Range("E20").Value = 3
$ws.Range("G20").Value = 66.31187066666666
$ws.Range("H20").Value = 198.935612
$ws.Range("I20").Value = 0.2937592662777732
$ws.Range("J20").Value = 0.2937592662777732
$ws.Range("K20").Value = 3
$ws.Range("M20").Value = 24.43393833333333
$ws.Range("N20").Value = 73.301815
$ws.Range("O20").Value = 0.1197493722400791
$ws.Range("P20").Value = 0.1197493722400791
$ws.Range("Q20").Value = 1620.260158637309
$ws.Range("R20").Value = 14582.34142773578
$ws.Range("S20").Value = 0.03517748772646958
$ws.Range("T20").Value = 0.03517748772646958

$ws.Range("E21").Value = 3
$ws.Range("G21").Value = 66.31187066666666
$ws.Range("H21").Value = 198.935612
$ws.Range("I21").Value = 0.2937592662777732
$ws.Range("J21").Value = 0.2937592662777732
$ws.Range("K21").Value = 3
$ws.Range("M21").Value = 86.43758666666668
$ws.Range("N21").Value = 259.31276
$ws.Range("O21").Value = 0.4236258027695808
$ws.Range("P21").Value = 0.4236258027695808
$ws.Range("Q21").Value = 5731.838067778792
$ws.Range("R21").Value = 51586.54261000912
$ws.Range("S21").Value = 0.1244440049979247
$ws.Range("T21").Value = 0.1244440049979247

$ws.Range("E22").Value = 3
$ws.Range("G22").Value = 66.31187066666666
$ws.Range("H22").Value = 198.935612
$ws.Range("I22").Value = 0.2937592662777732
$ws.Range("J22").Value = 0.2937592662777732
$ws.Range("K22").Value = 3
$ws.Range("M22").Value = 33.80250733333333
$ws.Range("N22").Value = 101.407522
$ws.Range("O22").Value = 0.1656642076314483
$ws.Range("P22").Value = 0.1656642076314483
$ws.Range("Q22").Value = 2241.507494497051
$ws.Range("R22").Value = 20173.56745047346
$ws.Range("S22").Value = 0.04866539608230292
$ws.Range("T22").Value = 0.04866539608230293

$ws.Range("E23").Value = 3
$ws.Range("G23").Value = 66.31187066666666
$ws.Range("H23").Value = 198.935612
$ws.Range("I23").Value = 0.2937592662777732
$ws.Range("J23").Value = 0.2937592662777732
$ws.Range("K23").Value = 3
$ws.Range("M23").Value = 25.01395866666667
$ws.Range("N23").Value = 75.041876
$ws.Range("O23").Value = 0.1225920196207674
$ws.Range("P23").Value = 0.1225920196207674
$ws.Range("Q23").Value = 1658.722391965346
$ws.Range("R23").Value = 14928.50152768811
$ws.Range("S23").Value = 0.03601254173530699
$ws.Range("T23").Value = 0.036012541735307

$ws.Range("E24").Value = 3
$ws.Range("G24").Value = 66.31187066666666
$ws.Range("H24").Value = 198.935612
$ws.Range("I24").Value = 0.2937592662777732
$ws.Range("J24").Value = 0.2937592662777732
$ws.Range("K24").Value = 3
$ws.Range("M24").Value = 25.369122
$ws.Range("N24").Value = 76.107366
$ws.Range("O24").Value = 0.1243326553557499
$ws.Range("P24").Value = 0.1243326553557499
$ws.Range("Q24").Value = 1682.273936990888
$ws.Range("R24").Value = 15140.46543291799
$ws.Range("S24").Value = 0.03652386961167235
$ws.Range("T24").Value = 0.03652386961167235

$ws.Range("E25").Value = 3
$ws.Range("G25").Value = 66.31187066666666
$ws.Range("H25").Value = 198.935612
$ws.Range("I25").Value = 0.2937592662777732
$ws.Range("J25").Value = 0.2937592662777732
$ws.Range("K25").Value = 3
$ws.Range("M25").Value = 8.985195333333333
$ws.Range("N25").Value = 26.955586
$ws.Range("O25").Value = 0.04403594238237437
$ws.Range("P25").Value = 0.04403594238237437
$ws.Range("Q25").Value = 595.8251108587368
$ws.Range("R25").Value = 5362.425997728632
$ws.Range("S25").Value = 0.01293596612409659
$ws.Range("T25").Value = 0.01293596612409659

$ws.Range("E26").Value = 3
$ws.Range("G26").Value = 59.79288466666666
$ws.Range("H26").Value = 179.378654
$ws.Range("I26").Value = 0.2648803864485285
$ws.Range("J26").Value = 0.2648803864485286
$ws.Range("K26").Value = 3
$ws.Range("M26").Value = 24.43393833333333
$ws.Range("N26").Value = 73.301815
$ws.Range("O26").Value = 0.1197493722400791
$ws.Range("P26").Value = 0.1197493722400791
$ws.Range("Q26").Value = 1460.975656717445
$ws.Range("R26").Value = 13148.78091045701
$ws.Range("S26").Value = 0.03171925999592085
$ws.Range("T26").Value = 0.03171925999592086

$ws.Range("E27").Value = 3
$ws.Range("G27").Value = 59.79288466666666
$ws.Range("H27").Value = 179.378654
$ws.Range("I27").Value = 0.2648803864485285
$ws.Range("J27").Value = 0.2648803864485286
$ws.Range("K27").Value = 3
$ws.Range("M27").Value = 86.43758666666668
$ws.Range("N27").Value = 259.31276
$ws.Range("O27").Value = 0.4236258027695808
$ws.Range("P27").Value = 0.4236258027695808
$ws.Range("Q27").Value = 5168.352650425004
$ws.Range("R27").Value = 46515.17385382504
$ws.Range("S27").Value = 0.1122101663471747
$ws.Range("T27").Value = 0.1122101663471747

$ws.Range("E28").Value = 3
$ws.Range("G28").Value = 59.79288466666666
$ws.Range("H28").Value = 179.378654
$ws.Range("I28").Value = 0.2648803864485285
$ws.Range("J28").Value = 0.2648803864485286
$ws.Range("K28").Value = 3
$ws.Range("M28").Value = 33.80250733333333
$ws.Range("N28").Value = 101.407522
$ws.Range("O28").Value = 0.1656642076314483
$ws.Range("P28").Value = 0.1656642076314483
$ws.Range("Q28").Value = 2021.149422426154
$ws.Range("R28").Value = 18190.34480183538
$ws.Range("S28").Value = 0.04388119933810729
$ws.Range("T28").Value = 0.04388119933810731

$ws.Range("E29").Value = 3
$ws.Range("G29").Value = 59.79288466666666
$ws.Range("H29").Value = 179.378654
$ws.Range("I29").Value = 0.2648803864485285
$ws.Range("J29").Value = 0.2648803864485286
$ws.Range("K29").Value = 3
$ws.Range("M29").Value = 25.01395866666667
$ws.Range("N29").Value = 75.041876
$ws.Range("O29").Value = 0.1225920196207674
$ws.Range("P29").Value = 0.1225920196207674
$ws.Range("Q29").Value = 1495.656745612767
$ws.Range("R29").Value = 13460.9107105149
$ws.Range("S29").Value = 0.03247222153265445
$ws.Range("T29").Value = 0.03247222153265446

$ws.Range("E30").Value = 3
$ws.Range("G30").Value = 59.79288466666666
$ws.Range("H30").Value = 179.378654
$ws.Range("I30").Value = 0.2648803864485285
$ws.Range("J30").Value = 0.2648803864485286
$ws.Range("K30").Value = 3
$ws.Range("M30").Value = 25.369122
$ws.Range("N30").Value = 76.107366
$ws.Range("O30").Value = 0.1243326553557499
$ws.Range("P30").Value = 0.1243326553557499
$ws.Range("Q30").Value = 1516.892985840596
$ws.Range("R30").Value = 13652.03687256536
$ws.Range("S30").Value = 0.03293328179880275
$ws.Range("T30").Value = 0.03293328179880276

$ws.Range("E31").Value = 3
$ws.Range("G31").Value = 59.79288466666666
$ws.Range("H31").Value = 179.378654
$ws.Range("I31").Value = 0.2648803864485285
$ws.Range("J31").Value = 0.2648803864485286
$ws.Range("K31").Value = 3
$ws.Range("M31").Value = 8.985195333333333
$ws.Range("N31").Value = 26.955586
$ws.Range("O31").Value = 0.04403594238237437
$ws.Range("P31").Value = 0.04403594238237437
$ws.Range("Q31").Value = 537.2507482734715
$ws.Range("R31").Value = 4835.256734461243
$ws.Range("S31").Value = 0.01166425743586846
$ws.Range("T31").Value = 0.01166425743586846

$ws.Range("E32").Value = 3
$ws.Range("G32").Value = 10.19562
$ws.Range("H32").Value = 30.58686
$ws.Range("I32").Value = 0.045166239774812
$ws.Range("J32").Value = 0.04516623977481201
$ws.Range("K32").Value = 3
$ws.Range("M32").Value = 24.43393833333333
$ws.Range("N32").Value = 73.301815
$ws.Range("O32").Value = 0.1197493722400791
$ws.Range("P32").Value = 0.1197493722400791
$ws.Range("Q32").Value = 249.1191503501
$ws.Range("R32").Value = 2242.0723531509
$ws.Range("S32").Value = 0.005408628859478629
$ws.Range("T32").Value = 0.005408628859478631

$ws.Range("E33").Value = 3
$ws.Range("G33").Value = 10.19562
$ws.Range("H33").Value = 30.58686
$ws.Range("I33").Value = 0.045166239774812
$ws.Range("J33").Value = 0.04516623977481201
$ws.Range("K33").Value = 3
$ws.Range("M33").Value = 86.43758666666668
$ws.Range("N33").Value = 259.31276
$ws.Range("O33").Value = 0.4236258027695808
$ws.Range("P33").Value = 0.4236258027695808
$ws.Range("Q33").Value = 881.2847873704001
$ws.Range("R33").Value = 7931.563086333601
$ws.Range("S33").Value = 0.01913358458268811
$ws.Range("T33").Value = 0.01913358458268811

$ws.Range("E34").Value = 3
$ws.Range("G34").Value = 10.19562
$ws.Range("H34").Value = 30.58686
$ws.Range("I34").Value = 0.045166239774812
$ws.Range("J34").Value = 0.04516623977481201
$ws.Range("K34").Value = 3
$ws.Range("M34").Value = 33.80250733333333
$ws.Range("N34").Value = 101.407522
$ws.Range("O34").Value = 0.1656642076314483
$ws.Range("P34").Value = 0.1656642076314483
$ws.Range("Q34").Value = 344.63751981788
$ws.Range("R34").Value = 3101.73767836092
$ws.Range("S34").Value = 0.007482429323986234
$ws.Range("T34").Value = 0.007482429323986237

$ws.Range("E35").Value = 3
$ws.Range("G35").Value = 10.19562
$ws.Range("H35").Value = 30.58686
$ws.Range("I35").Value = 0.045166239774812
$ws.Range("J35").Value = 0.04516623977481201
$ws.Range("K35").Value = 3
$ws.Range("M35").Value = 25.01395866666667
$ws.Range("N35").Value = 75.041876
$ws.Range("O35").Value = 0.1225920196207674
$ws.Range("P35").Value = 0.1225920196207674
$ws.Range("Q35").Value = 255.03281726104
$ws.Range("R35").Value = 2295.29535534936
$ws.Range("S35").Value = 0.005537020552670037
$ws.Range("T35").Value = 0.005537020552670038

$ws.Range("E36").Value = 3
$ws.Range("G36").Value = 10.19562
$ws.Range("H36").Value = 30.58686
$ws.Range("I36").Value = 0.045166239774812
$ws.Range("J36").Value = 0.04516623977481201
$ws.Range("K36").Value = 3
$ws.Range("M36").Value = 25.369122
$ws.Range("N36").Value = 76.107366
$ws.Range("O36").Value = 0.1243326553557499
$ws.Range("P36").Value = 0.1243326553557499
$ws.Range("Q36").Value = 258.65392764564
$ws.Range("R36").Value = 2327.88534881076
$ws.Range("S36").Value = 0.005615638523636866
$ws.Range("T36").Value = 0.005615638523636867

$ws.Range("E37").Value = 3
$ws.Range("G37").Value = 10.19562
$ws.Range("H37").Value = 30.58686
$ws.Range("I37").Value = 0.045166239774812
$ws.Range("J37").Value = 0.04516623977481201
$ws.Range("K37").Value = 3
$ws.Range("M37").Value = 8.985195333333333
$ws.Range("N37").Value = 26.955586
$ws.Range("O37").Value = 0.04403594238237437
$ws.Range("P37").Value = 0.04403594238237437
$ws.Range("Q37").Value = 91.60963724444
$ws.Range("R37").Value = 824.4867351999601
$ws.Range("S37").Value = 0.001988937932352127
$ws.Range("T37").Value = 0.001988937932352127
